$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting from existing header cell H1 onto the two new header cells
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# Values for column I (I0) and column J (IF), rows 2-55
$iValues = @(6,5,7,4,5,5,3,9,6,7,7,9,7,8,9,9,6,4,6,3,4,8,8,9,6,8,9,6,8,1,5,5,5,7,6,8,7,7,7,7,9,6,7,7,6,5,6,7,6,5,4,6,5,6)
$jValues = @(7,6,8,6,5,6,5,9,7,8,7,9,9,8,9,9,8,6,6,6,6,8,9,9,6,8,9,8,8,3,6,7,6,8,6,9,7,8,8,8,9,7,9,8,8,6,8,8,8,6,5,6,6,7)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
